{"js": "// Replace the multiplication-fact answers throughout the document's\n// table cells with the newly generated values. Each old value is\n// unique in the document, so a scoped search-and-replace on the\n// document body is unambiguous and preserves run formatting\n// (font/size) because we replace only the text of the matched range.\nconst replacements = [\n  [\"21\u00d724=504\", \"75\u00d754=4050\"],\n  [\"43\u00d769=2967\", \"41\u00d723=943\"],\n  [\"47\u00d794=4418\", \"35\u00d745=1575\"],\n  [\"19\u00d722=418\", \"24\u00d751=1224\"],\n  [\"89\u00d746=4094\", \"57\u00d753=3021\"],\n  [\"71\u00d754=3834\", \"43\u00d774=3182\"],\n  [\"88\u00d774=6512\", \"32\u00d793=2976\"],\n  [\"46\u00d753=2438\", \"53\u00d751=2703\"],\n  [\"83\u00d799=8217\", \"19\u00d779=1501\"],\n  [\"76\u00d795=7220\", \"37\u00d755=2035\"],\n  [\"37\u00d743=1591\", \"31\u00d771=2201\"],\n  [\"58\u00d752=3016\", \"96\u00d735=3360\"],\n  [\"44\u00d794=4136\", \"78\u00d724=1872\"],\n  [\"15\u00d749=735\", \"79\u00d739=3081\"],\n  [\"88\u00d783=7304\", \"56\u00d789=4984\"],\n  [\"78\u00d733=2574\", \"35\u00d731=1085\"],\n  [\"70\u00d764=4480\", \"12\u00d750=600\"],\n  [\"96\u00d787=8352\", \"34\u00d736=1224\"],\n  [\"54\u00d732=1728\", \"94\u00d786=8084\"],\n  [\"68\u00d731=2108\", \"96\u00d766=6336\"],\n  [\"65\u00d759=3835\", \"47\u00d777=3619\"],\n  [\"41\u00d747=1927\", \"55\u00d758=3190\"],\n  [\"94\u00d771=6674\", \"67\u00d721=1407\"],\n  [\"30\u00d719=570\", \"15\u00d770=1050\"],\n  [\"41\u00d753=2173\", \"31\u00d789=2759\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the multiplication-fact answers throughout the document's\n# table cells with the newly generated values. Each old value is\n# unique in the document, so a document-wide Find/Replace per pair is\n# unambiguous and (with no formatting specified on Find/Replacement)\n# preserves each run's existing formatting (font/size).\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"21\u00d724=504\", \"75\u00d754=4050\"),\n  @(\"43\u00d769=2967\", \"41\u00d723=943\"),\n  @(\"47\u00d794=4418\", \"35\u00d745=1575\"),\n  @(\"19\u00d722=418\", \"24\u00d751=1224\"),\n  @(\"89\u00d746=4094\", \"57\u00d753=3021\"),\n  @(\"71\u00d754=3834\", \"43\u00d774=3182\"),\n  @(\"88\u00d774=6512\", \"32\u00d793=2976\"),\n  @(\"46\u00d753=2438\", \"53\u00d751=2703\"),\n  @(\"83\u00d799=8217\", \"19\u00d779=1501\"),\n  @(\"76\u00d795=7220\", \"37\u00d755=2035\"),\n  @(\"37\u00d743=1591\", \"31\u00d771=2201\"),\n  @(\"58\u00d752=3016\", \"96\u00d735=3360\"),\n  @(\"44\u00d794=4136\", \"78\u00d724=1872\"),\n  @(\"15\u00d749=735\", \"79\u00d739=3081\"),\n  @(\"88\u00d783=7304\", \"56\u00d789=4984\"),\n  @(\"78\u00d733=2574\", \"35\u00d731=1085\"),\n  @(\"70\u00d764=4480\", \"12\u00d750=600\"),\n  @(\"96\u00d787=8352\", \"34\u00d736=1224\"),\n  @(\"54\u00d732=1728\", \"94\u00d786=8084\"),\n  @(\"68\u00d731=2108\", \"96\u00d766=6336\"),\n  @(\"65\u00d759=3835\", \"47\u00d777=3619\"),\n  @(\"41\u00d747=1927\", \"55\u00d758=3190\"),\n  @(\"94\u00d771=6674\", \"67\u00d721=1407\"),\n  @(\"30\u00d719=570\", \"15\u00d770=1050\"),\n  @(\"41\u00d753=2173\", \"31\u00d789=2759\")\n)\n\nforeach ($pair in $pairs) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}\n"}
